# "small fix in all templates"
#
# Changes applied:
#  - Rename the first sheet from "SampleSubmission" to "SampleUpdate"
#  - Make the first sheet ("SampleUpdate") the active/selected tab
#    (previously "Index" was the active tab)
#  - Update the selection on "SampleUpdate" from E8 to B8
#  - "Index" sheet is no longer the active tab

$wb = $excel.ActiveWorkbook

$wsUpdate = $wb.Worksheets.Item(1)
$wsIndex  = $wb.Worksheets.Item(2)

# Rename SampleSubmission -> SampleUpdate
$wsUpdate.Name = "SampleUpdate"

# Leave the Index sheet's own selection untouched, but switch the active
# tab to SampleUpdate and select B8 there.
$wsIndex.Activate() | Out-Null
$wsUpdate.Activate() | Out-Null
$wsUpdate.Range("B8").Select() | Out-Null

Write-Output "Renamed sheet 1 to '$($wsUpdate.Name)' and selected B8 as the active cell on the active tab."
